$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# Two new rows were appended to the field-list table (SEQ 35 & 36).
# Copy the formatting from the row above first (format-only paste) so the
# newly-populated cells keep the existing column look (style 15/16) instead
# of falling back to the column's default style.
$ws.Range("A42:E42").Copy()
$ws.Range("A43:E43").PasteSpecial(-4122)
$ws.Range("A42:E42").Copy()
$ws.Range("A44:E44").PasteSpecial(-4122)

$ws.Range("A43").Value = 35
$ws.Range("B43").Value = "ActualFilingDate"
$ws.Range("C43").Value = "實際報送日期"
$ws.Range("D43").Value = "Decimald"
$ws.Range("E43").Value = 8

$ws.Range("A44").Value = 36
$ws.Range("B44").Value = "ActualFilingMark"
$ws.Range("C44").Value = "實際報送記號"
$ws.Range("D44").Value = "VARCHAR2"
$ws.Range("E44").Value = 3

# The cursor ended up on D48 when the edit was committed.
$ws.Range("D48").Select()
